# Max path binary sum
# Adds Notes / Time Complexity / Space Complexity for the
# "Binary Tree Maximum Path Sum" row (row 8) and updates the
# sheet view selection / scroll position accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the new cell content for row 8 (D8:F8)
$ws.Range("D8").Value = "#Post traversal`nDo a dfs`nKeep max_sum which should be the sum of the left max path and the right max path + the node val`n#remember when you return you will return only the node.val and max of left_sum or right_sum. reason is as you pass on the value to the root you can only pass left or right max path"

$ws.Range("E8").Value = "dfs visits each node and hence time complexity is O(n)"

$ws.Range("F8").Value = " The main factor contributing to the space complexity (beyond the storage for the tree itself) is the recursion call stack used by the dfs function. The maximum depth of the recursion call stack corresponds to the height (H) of the binary tree. . Worst case it is a skewed "

# Match styling used by the rest of the sheet (wrap text) and row height
$ws.Range("D8:F8").WrapText = $true
$ws.Rows.Item(8).RowHeight = 101.5

# Update the view: scrolled/selected cell moved from D8 to D5
$ws.Range("D5").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 3
